$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while preserving it as text (avoids Excel
# auto-converting numeric-looking strings like "1.00" into numbers),
# and without leaving a residual NumberFormat/style on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# --- Row 2..51 price/volume updates ---
Set-TextValue $ws.Range("D2") "67.334.28"
Set-TextValue $ws.Range("E2") "  +1.63%  "
Set-TextValue $ws.Range("D3") "3.879.34"
Set-TextValue $ws.Range("E3") "  +1.22%  "
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  -0.06%  "
Set-TextValue $ws.Range("D5") "471.32"
Set-TextValue $ws.Range("E5") "  +10.11%  "
Set-TextValue $ws.Range("D6") "145.72"
Set-TextValue $ws.Range("E6") "  +10.52%  "
Set-TextValue $ws.Range("D7") "0.635"
Set-TextValue $ws.Range("E7") "  +3.35%  "
Set-TextValue $ws.Range("D8") "0.999"
Set-TextValue $ws.Range("E8") "  -0.14%  "
Set-TextValue $ws.Range("D9") "0.748"
Set-TextValue $ws.Range("E9") "  +1.91%  "
Set-TextValue $ws.Range("D10") "0.155"
Set-TextValue $ws.Range("E10") "  -1.55%  "
Set-TextValue $ws.Range("E11") "  -6.27%  "
Set-TextValue $ws.Range("D12") "43.53"
Set-TextValue $ws.Range("E12") "  +3.96%  "
Set-TextValue $ws.Range("D13") "10.47"
Set-TextValue $ws.Range("E13") "  +0.12%  "
Set-TextValue $ws.Range("D14") "4.505.35"
Set-TextValue $ws.Range("E14") "  +1.33%  "
Set-TextValue $ws.Range("D15") "14.85"
Set-TextValue $ws.Range("E15") "  -5.83%  "
Set-TextValue $ws.Range("D16") "3.860.17"
Set-TextValue $ws.Range("E16") "  +1.09%  "
Set-TextValue $ws.Range("E17") "  -0.47%  "
Set-TextValue $ws.Range("D18") "20.12"
Set-TextValue $ws.Range("E18") "  +0.11%  "
Set-TextValue $ws.Range("E19") "  +6.16%  "
Set-TextValue $ws.Range("D20") "67.627.77"
Set-TextValue $ws.Range("E20") "  +1.64%  "
Set-TextValue $ws.Range("D21") "437.39"
Set-TextValue $ws.Range("E21") "  +5.24%  "
Set-TextValue $ws.Range("E22") "  -1.24%  "
Set-TextValue $ws.Range("E23") "  +6.01%  "
Set-TextValue $ws.Range("D24") "89.44"
Set-TextValue $ws.Range("E24") "  +4.96%  "
Set-TextValue $ws.Range("D25") "3.62"
Set-TextValue $ws.Range("E25") "  +10.02%  "
Set-TextValue $ws.Range("D26") "38.13"
Set-TextValue $ws.Range("E26") "  +1.88%  "
Set-TextValue $ws.Range("D27") "10.09"
Set-TextValue $ws.Range("E27") "  +6.78%  "
Set-TextValue $ws.Range("D28") "9.97"
Set-TextValue $ws.Range("E28") "  -1.45%  "
Set-TextValue $ws.Range("D29") "5.49"
Set-TextValue $ws.Range("E29") "  +2.33%  "
Set-TextValue $ws.Range("D30") "732.53"
Set-TextValue $ws.Range("E30") "  +1.54%  "
Set-TextValue $ws.Range("D31") "13.89"
Set-TextValue $ws.Range("E31") "  +0.04%  "
Set-TextValue $ws.Range("D32") "0.134"
Set-TextValue $ws.Range("E32") "  +6.58%  "
Set-TextValue $ws.Range("D33") "2.78"
Set-TextValue $ws.Range("E33") "  +0.26%  "
Set-TextValue $ws.Range("D34") "44.45"
Set-TextValue $ws.Range("E34") "  +12.99%  "
Set-TextValue $ws.Range("D35") "0.163"
Set-TextValue $ws.Range("E35") "  +8.09%  "
Set-TextValue $ws.Range("D36") "58.08"
Set-TextValue $ws.Range("E36") "  +4.38%  "
Set-TextValue $ws.Range("E37") "  +0.19%  "
Set-TextValue $ws.Range("D38") "5.53"
Set-TextValue $ws.Range("E38") "  -4.86%  "
Set-TextValue $ws.Range("D39") "0.0486"
Set-TextValue $ws.Range("E39") "  +3.10%  "
Set-TextValue $ws.Range("D40") "0.349"
Set-TextValue $ws.Range("E40") "  +7.53%  "
Set-TextValue $ws.Range("D41") "2.93"
Set-TextValue $ws.Range("E41") "  +0.97%  "
Set-TextValue $ws.Range("D42") "0.0₃0688"
Set-TextValue $ws.Range("E42") "  -6.03%  "
Set-TextValue $ws.Range("E43") "  +3.30%  "
Set-TextValue $ws.Range("E44") "  +0.13%  "
Set-TextValue $ws.Range("E45") "  +4.03%  "
Set-TextValue $ws.Range("E46") "  +1.90%  "
Set-TextValue $ws.Range("D47") "3.27"
Set-TextValue $ws.Range("E47") "  -0.39%  "
Set-TextValue $ws.Range("D50") "2.90"
Set-TextValue $ws.Range("E50") "  +0.83%  "
Set-TextValue $ws.Range("D51") "144.31"
Set-TextValue $ws.Range("E51") "  +1.51%  "

# --- Rows 48 and 49: coin identities swapped (ARBITRUM <-> WEMIXToken) ---
Set-TextValue $ws.Range("B48") "WEMIXToken"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D48") "2.76"
Set-TextValue $ws.Range("E48") "  +4.70%  "

Set-TextValue $ws.Range("B49") "ARBITRUM"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D49") "2.17"
Set-TextValue $ws.Range("E49") "  +4.87%  "
